# Update cryptos price/volume data to reflect the latest scrape.
# Generated from the authoritative cell-level diff; each row below lists
# only the columns that actually changed for that coin.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value to a cell while always preserving it as TEXT
# (the sheet keeps prices/percentages as strings, e.g. "0.9996" or "1.000",
# and Excel would otherwise silently coerce plain-numeric-looking strings
# into real numbers).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "22.365.50"
Set-TextValue $ws.Range("E2") "  -1.01%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.562.44"
Set-TextValue $ws.Range("E3") "  -1.14%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.9996"
Set-TextValue $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "0.9999"
Set-TextValue $ws.Range("E5") "  -0.03%  "

# Row 6
Set-TextValue $ws.Range("D6") "287.92"
Set-TextValue $ws.Range("E6") "  -0.29%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +3.22%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3284"
Set-TextValue $ws.Range("E8") "  -2.55%  "

# Row 9
Set-TextValue $ws.Range("D9") "44.52"
Set-TextValue $ws.Range("E9") "  -8.45%  "

# Row 10
Set-TextValue $ws.Range("D10") "1.144"
Set-TextValue $ws.Range("E10") "  +0.03%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.07395"
Set-TextValue $ws.Range("E11") "  -1.38%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.9996"
Set-TextValue $ws.Range("E12") "  -0.01%  "

# Row 13
Set-TextValue $ws.Range("D13") "20.49"
Set-TextValue $ws.Range("E13") "  -3.08%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.873"
Set-TextValue $ws.Range("E14") "  -2.70%  "

# Row 15
Set-TextValue $ws.Range("D15") "6.816"
Set-TextValue $ws.Range("E15") "  -2.25%  "

# Row 16
Set-TextValue $ws.Range("D16") "1.559.24"
Set-TextValue $ws.Range("E16") "  -0.71%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.00001088"
Set-TextValue $ws.Range("E17") "  -3.00%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.06685"
Set-TextValue $ws.Range("E18") "  -1.32%  "

# Row 19
Set-TextValue $ws.Range("D19") "86.39"
Set-TextValue $ws.Range("E19") "  -2.88%  "

# Row 20
Set-TextValue $ws.Range("D20") "6.458"
Set-TextValue $ws.Range("E20") "  +0.51%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.9998"
Set-TextValue $ws.Range("E21") "  -0.01%  "

# Row 22
Set-TextValue $ws.Range("D22") "16.25"
Set-TextValue $ws.Range("E22") "  -2.17%  "

# Row 23
Set-TextValue $ws.Range("D23") "11.71"
Set-TextValue $ws.Range("E23") "  -4.03%  "

# Row 24
Set-TextValue $ws.Range("D24") "22.358.25"
Set-TextValue $ws.Range("E24") "  -1.01%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.294"
Set-TextValue $ws.Range("E25") "  -3.87%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.594"
Set-TextValue $ws.Range("E26") "  -1.60%  "

# Row 27
Set-TextValue $ws.Range("D27") "150.00"
Set-TextValue $ws.Range("E27") "  -2.00%  "

# Row 28
Set-TextValue $ws.Range("D28") "19.44"
Set-TextValue $ws.Range("E28") "  -1.30%  "

# Row 29
Set-TextValue $ws.Range("D29") "4.936"
Set-TextValue $ws.Range("E29") "  -1.72%  "

# Row 30
Set-TextValue $ws.Range("D30") "122.73"
Set-TextValue $ws.Range("E30") "  -1.60%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.740.08"
Set-TextValue $ws.Range("E31") "  -0.74%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.083"
Set-TextValue $ws.Range("E32") "  +0.41%  "

# Row 33
Set-TextValue $ws.Range("D33") "5.973"
Set-TextValue $ws.Range("E33") "  -4.07%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.892"
Set-TextValue $ws.Range("E34") "  -5.53%  "

# Row 35
Set-TextValue $ws.Range("D35") "9.462"
Set-TextValue $ws.Range("E35") "  -3.97%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.08320"
Set-TextValue $ws.Range("E36") "  -0.15%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02395"
Set-TextValue $ws.Range("E37") "  -3.33%  "

# Row 38
Set-TextValue $ws.Range("D38") "5.353"
Set-TextValue $ws.Range("E38") "  -2.61%  "

# Row 39
Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.06302"
Set-TextValue $ws.Range("E39") "  -2.66%  "

# Row 40
Set-TextValue $ws.Range("B40") "TrustWalletToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D40") "1.280"
Set-TextValue $ws.Range("E40") "  -2.05%  "

# Row 41
Set-TextValue $ws.Range("B41") "Algorand"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D41") "0.2176"
Set-TextValue $ws.Range("E41") "  -4.60%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -2.91%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.6104"
Set-TextValue $ws.Range("E43") "  -3.90%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.000"
Set-TextValue $ws.Range("E44") "  -0.01%  "

# Row 45
Set-TextValue $ws.Range("D45") "13.81"
Set-TextValue $ws.Range("E45") "  -2.35%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.5952"
Set-TextValue $ws.Range("E46") "  -4.01%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.750"
Set-TextValue $ws.Range("E47") "  -0.70%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.008"
Set-TextValue $ws.Range("E48") "  -3.35%  "

# Row 49
Set-TextValue $ws.Range("E49") "  -1.31%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.180"

# Row 51
Set-TextValue $ws.Range("D51") "0.07108"
